$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.076.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.010.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.85%  "
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.607"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.96%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.12"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.19%  "
$ws.Range("E9").Value = "  -2.97%  "
$ws.Range("E10").Value = "  -4.67%  "
$ws.Range("E11").Value = "  -4.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.305.89"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.81%  "
$ws.Range("E13").Value = "  -3.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.67%  "
$ws.Range("E16").Value = "  -2.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.011.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.960.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("E19").Value = "  +2.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0810"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.13%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  +1.95%  "
$ws.Range("E25").Value = "  -5.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.71%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.124"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.80%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.116"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0600"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.45"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.55%  "
$ws.Range("E35").Value = "  -2.69%  "
$ws.Range("E36").Value = "  +2.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("E38").Value = "  -2.51%  "
$ws.Range("E39").Value = "  -1.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.458.55"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0212"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "94.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0912"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.69%  "
$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.19"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +15.55%  "
$ws.Range("B45").Value = "HuobiToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.76"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.23%  "
$ws.Range("E46").Value = "  -2.64%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.94%  "
$ws.Range("E48").Value = "  -2.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.05"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.55%  "
$ws.Range("E50").Value = "  -0.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.196.40"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.66%  "
